# Automatic update of files.
# Rows were re-matched to different source records, so the editable fields
# of each row pair got swapped between the two rows.
# Pairs: (10,11) (14,15) (16,17) (20,21)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between the two rows of each swapped pair.
# (Columns that already hold identical values on both rows, e.g. the
# species/location/date columns that didn't change, are left untouched.)
$colsCommon = @("A", "Q", "R", "X", "Z", "AB", "AX")

function Swap-Cell {
    param($col, $row1, $row2)

    $c1 = $ws.Range($col + $row1)
    $c2 = $ws.Range($col + $row2)

    $v1 = $c1.Value2
    $v2 = $c2.Value2

    $c1.Value = $v2
    $c2.Value = $v1
}

function Swap-TextCell {
    # Like Swap-Cell, but forces the cell to keep/become a TEXT value even
    # when the content looks numeric (e.g. "1", "2"), matching the column's
    # existing storage convention in the sheet.
    param($col, $row1, $row2)

    $c1 = $ws.Range($col + $row1)
    $c2 = $ws.Range($col + $row2)

    $v1 = $c1.Text
    $v2 = $c2.Text

    $c1.NumberFormat = "@"
    $c2.NumberFormat = "@"

    if ($v2 -eq "") {
        $c1.Value = ""
    } else {
        $c1.Value = $v2
    }

    if ($v1 -eq "") {
        $c2.Value = ""
    } else {
        $c2.Value = $v1
    }
}

function Swap-Row {
    param($row1, $row2, $includeJ)

    foreach ($col in $colsCommon) {
        Swap-Cell $col $row1 $row2
    }

    Swap-TextCell "I" $row1 $row2

    if ($includeJ) {
        Swap-Cell "J" $row1 $row2
    }
}

# Rows 10 / 11 : only the common + I columns differ
Swap-Row 10 11 $false

# Rows 14 / 15 : only the common + I columns differ
Swap-Row 14 15 $false

# Rows 16 / 17 : whole record swapped (different species / location too)
foreach ($col in @("A","B","E","F","G","H","P","Q","R","X","Z","AB","AX")) {
    Swap-Cell $col 16 17
}

# Rows 20 / 21 : common + I + J columns differ
Swap-Row 20 21 $true
